$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextCell "D2" "26.866.24"
Set-TextCell "E2" "  +0.37%  "
Set-TextCell "D3" "1.642.67"
Set-TextCell "E3" "  -0.06%  "
Set-TextCell "D5" "218.66"
Set-TextCell "E5" "  +0.91%  "
Set-TextCell "D6" "0.498"
Set-TextCell "E6" "  -0.28%  "
Set-TextCell "E7" "  -0.14%  "
Set-TextCell "E8" "  +0.02%  "
Set-TextCell "E9" "  -1.12%  "
Set-TextCell "E10" "  +0.59%  "
Set-TextCell "D11" "0.0846"
Set-TextCell "E11" "  +0.50%  "
Set-TextCell "D12" "1.871.77"
Set-TextCell "E12" "  +0.09%  "
Set-TextCell "D13" "1.632.25"
Set-TextCell "E13" "  -0.73%  "
Set-TextCell "D14" "4.17"
Set-TextCell "E14" "  +0.17%  "
Set-TextCell "E15" "  +0.11%  "
Set-TextCell "D16" "65.49"
Set-TextCell "E16" "  +1.93%  "
Set-TextCell "D17" "26.869.34"
Set-TextCell "E17" "  +0.43%  "
Set-TextCell "D18" "0.0₃0734"
Set-TextCell "E18" "  -0.49%  "
Set-TextCell "D19" "215.39"
Set-TextCell "E19" "  +0.71%  "
Set-TextCell "E20" "  -0.24%  "
Set-TextCell "D21" "6.69"
Set-TextCell "E21" "  +6.74%  "
Set-TextCell "E22" "  +0.24%  "
Set-TextCell "D23" "2.38"
Set-TextCell "E23" "  -1.14%  "
Set-TextCell "E24" "  -1.24%  "
Set-TextCell "D25" "147.58"
Set-TextCell "E25" "  +1.86%  "
Set-TextCell "E26" "  +0.13%  "
Set-TextCell "D27" "0.118"
Set-TextCell "E27" "  -0.09%  "
Set-TextCell "D28" "7.22"
Set-TextCell "E28" "  +1.82%  "
Set-TextCell "D29" "15.74"
Set-TextCell "E29" "  +0.61%  "
Set-TextCell "E30" "  -0.13%  "
Set-TextCell "E31" "  +1.19%  "
Set-TextCell "E32" "  +1.69%  "
Set-TextCell "E33" "  -0.08%  "
Set-TextCell "D34" "1.279.35"
Set-TextCell "E34" "  -1.13%  "
Set-TextCell "D35" "1.54"
Set-TextCell "E35" "  +0.52%  "
Set-TextCell "E36" "  -0.02%  "
Set-TextCell "E37" "  +0.53%  "
Set-TextCell "E38" "  +0.02%  "
Set-TextCell "D39" "0.818"
Set-TextCell "E39" "  -1.03%  "
Set-TextCell "E40" "  -0.13%  "
Set-TextCell "E41" "  -0.13%  "
Set-TextCell "D42" "5.35"
Set-TextCell "E42" "  +0.14%  "
Set-TextCell "D43" "1.782.84"
Set-TextCell "E43" "  -0.50%  "
Set-TextCell "D44" "2.12"
Set-TextCell "E44" "  -5.30%  "
Set-TextCell "D45" "92.80"
Set-TextCell "E45" "  +1.56%  "
Set-TextCell "D46" "61.06"
Set-TextCell "E46" "  -0.11%  "
Set-TextCell "E48" "  -0.28%  "
Set-TextCell "B49" "Algorand"
Set-TextCell "C49" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D49" "0.0966"
Set-TextCell "E49" "  -0.95%  "
Set-TextCell "B50" "EnergySwap"
Set-TextCell "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D50" "7.56"
Set-TextCell "E50" "  -1.15%  "
Set-TextCell "E51" "  -0.16%  "
